# aggiornamento fino a 27/05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44330, 10, 43, 106.7580316798252),
    @(44331, 1, 33, 81.93058245195888),
    @(44332, 4, 34, 84.41332737474552),
    @(44333, 6, 31, 76.96509260638562),
    @(44334, 1, 28, 69.51685783802573),
    @(44335, 3, 31, 76.96509260638562),
    @(44336, 8, 33, 81.93058245195888),
    @(44337, 3, 26, 64.55136799245246),
    @(44338, 3, 28, 69.51685783802573),
    @(44339, 2, 26, 64.55136799245246),
    @(44340, 1, 21, 52.13764337851929),
    @(44341, 0, 20, 49.65489845573266),
    @(44342, 1, 18, 44.68940861015939),
    @(44343, 1, 11, 27.31019415065296)
)

# Template cell whose formatting (date number format, border, centered
# alignment) must be propagated to the newly appended date cells in col A.
$ws.Range("A255").Copy()

$startRow = 256
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]

    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $rowVals[0]

    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
}
